$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account-statement rows (doc id, name, period, value) replacing the
# previous periods table in rows 16-34. Order comes from the updated
# source database: Ricardo Andres Barrios Montes' periods now run from
# 2306 down to 2202, followed by Natalia Andrea Ricardo Mena (2202) and
# Jose Daniel Villa Villarreal (2204) at the bottom.
$rows = @(
    @{ Row = 16; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2306"; Value = 58667 },
    @{ Row = 17; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2305"; Value = 80000 },
    @{ Row = 18; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2304"; Value = 80000 },
    @{ Row = 19; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2303"; Value = 80000 },
    @{ Row = 20; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2302"; Value = 80000 },
    @{ Row = 21; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2301"; Value = 80000 },
    @{ Row = 22; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2212"; Value = 80000 },
    @{ Row = 23; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2211"; Value = 80000 },
    @{ Row = 24; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2210"; Value = 80000 },
    @{ Row = 25; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2209"; Value = 80000 },
    @{ Row = 26; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2208"; Value = 80000 },
    @{ Row = 27; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2207"; Value = 80000 },
    @{ Row = 28; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2206"; Value = 80000 },
    @{ Row = 29; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2205"; Value = 80000 },
    @{ Row = 30; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2204"; Value = 80000 },
    @{ Row = 31; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2203"; Value = 80000 },
    @{ Row = 32; Doc = "1143384213"; Name = "RICARDO ANDRES BARRIOS MONTES"; Period = "2202"; Value = 80000 },
    @{ Row = 33; Doc = "1007980980"; Name = "NATALIA ANDREA RICARDO MENA";   Period = "2202"; Value = 80000 },
    @{ Row = 34; Doc = "1007229509"; Name = "JOSE DANIEL VILLA VILLARREAL"; Period = "2204"; Value = 53334 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc      # C: N° Doc Trabajador
    $ws.Cells.Item($r.Row, 4).Value = $r.Name     # D: Nombre Trabajador
    $ws.Cells.Item($r.Row, 5).Value = $r.Period   # E: Periodo Mora
    $ws.Cells.Item($r.Row, 6).Value = $r.Value    # F: Valor Mora
}
